# Updated remaining queries for C3DC
#
# The workbook stores the same DuckDB/SQL JOIN pattern in several cells
# (C2, B2, B3, B4, B5, B6, B7). Each occurrence of the old join keys
# (std.id / prt.id) is replaced with the new, explicit join keys
# (std.study_id / prt.participant_id), matching on the qualified
# "table.column" aliases too.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cells = @("C2", "B2", "B3", "B4", "B5", "B6", "B7")

foreach ($addr in $cells) {
    $range = $ws.Range($addr)
    $text = $range.Value2

    $text = $text -replace 'std\.id = prt\."study\.id"', 'std.study_id = prt."study.study_id"'
    $text = $text -replace 'prt\.id = dgn\."participant\.id"', 'prt.participant_id = dgn."participant.participant_id"'
    $text = $text -replace 'prt\.id = trt\."participant\.id"', 'prt.participant_id = trt."participant.participant_id"'
    $text = $text -replace 'prt\.id = trr\."participant\.id"', 'prt.participant_id = trr."participant.participant_id"'
    $text = $text -replace 'prt\.id = srv\."participant\.id"', 'prt.participant_id = srv."participant.participant_id"'
    $text = $text -replace 'std\.id = rfs\."study\.id"', 'std.study_id = rfs."study.study_id"'

    $range.Value = $text
}

# Column C was widened to fit the longer query text (and is no longer
# flagged as an auto "best fit" column).
$ws.Columns.Item(3).ColumnWidth = 67.83
